$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "ceti"
$ws.Range("B5").Value = "3499320013789"

$ws.Range("B8").Select()
